# Update "想去人数" (number of people interested) values that changed
# between the two data snapshots, as described by the commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 190
$ws1.Range("F3").Value = 236
$ws1.Range("F4").Value = 260
$ws1.Range("F5").Value = 792
$ws1.Range("F6").Value = 242
$ws1.Range("F7").Value = 5868
$ws1.Range("F8").Value = 33
$ws1.Range("F13").Value = 129
$ws1.Range("F15").Value = 351
$ws1.Range("F16").Value = 27

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 4

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 190
$ws4.Range("F3").Value = 236
$ws4.Range("F4").Value = 260
$ws4.Range("F5").Value = 792
$ws4.Range("F6").Value = 242
$ws4.Range("F7").Value = 5868
$ws4.Range("F8").Value = 33
$ws4.Range("F13").Value = 129
$ws4.Range("F15").Value = 351
$ws4.Range("F16").Value = 27
$ws4.Range("F17").Value = 4
